# Natmi following Dr Hou advice
# Recompute the Tnc-Egfr ligand-receptor pair table: refresh the existing
# Sending-cluster -> Target-cluster combinations for rows 2-7 with updated
# NATMI statistics, and append the 3 new combinations (rows 8-10) introduced
# by treating "ECs" as a distinct sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate, G=Ligand average expression value,
# H=Ligand total expression value, I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate, M=Receptor average expression value,
# N=Receptor total expression value, O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity

# Row 2: ECs -> ECs  (Tnc/Egfr)
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tnc"
$ws.Cells.Item(2,3).Value = "Egfr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.442371333333333
$ws.Cells.Item(2,8).Value = 4.327114
$ws.Cells.Item(2,9).Value = 0.03522044016446201
$ws.Cells.Item(2,10).Value = 0.03522044016446201
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.307106666666667
$ws.Cells.Item(2,14).Value = 3.92132
$ws.Cells.Item(2,15).Value = 0.01256263154946851
$ws.Cells.Item(2,16).Value = 0.01256263154946851
$ws.Cells.Item(2,17).Value = 1.885333185608889
$ws.Cells.Item(2,18).Value = 16.96799867048
$ws.Cells.Item(2,19).Value = 0.0004424614127962383
$ws.Cells.Item(2,20).Value = 0.0004424614127962383

# Row 3: ECs -> FAPs  (Tnc/Egfr)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tnc"
$ws.Cells.Item(3,3).Value = "Egfr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.442371333333333
$ws.Cells.Item(3,8).Value = 4.327114
$ws.Cells.Item(3,9).Value = 0.03522044016446201
$ws.Cells.Item(3,10).Value = 0.03522044016446201
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 80.22623699999998
$ws.Cells.Item(3,14).Value = 240.678711
$ws.Cells.Item(3,15).Value = 0.77105616682495
$ws.Cells.Item(3,16).Value = 0.77105616682495
$ws.Cells.Item(3,17).Value = 115.716024430006
$ws.Cells.Item(3,18).Value = 1041.444219870054
$ws.Cells.Item(3,19).Value = 0.02715693758709759
$ws.Cells.Item(3,20).Value = 0.02715693758709759

# Row 4: ECs -> sCs  (Tnc/Egfr)
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tnc"
$ws.Cells.Item(4,3).Value = "Egfr"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.442371333333333
$ws.Cells.Item(4,8).Value = 4.327114
$ws.Cells.Item(4,9).Value = 0.03522044016446201
$ws.Cells.Item(4,10).Value = 0.03522044016446201
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 22.51385866666667
$ws.Cells.Item(4,14).Value = 67.541576
$ws.Cells.Item(4,15).Value = 0.2163812016255815
$ws.Cells.Item(4,16).Value = 0.2163812016255815
$ws.Cells.Item(4,17).Value = 32.47334434351822
$ws.Cells.Item(4,18).Value = 292.260099091664
$ws.Cells.Item(4,19).Value = 0.007621041164568182
$ws.Cells.Item(4,20).Value = 0.007621041164568182

# Row 5: FAPs -> ECs  (Tnc/Egfr)
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Tnc"
$ws.Cells.Item(5,3).Value = "Egfr"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 11.331397
$ws.Cells.Item(5,8).Value = 33.994191
$ws.Cells.Item(5,9).Value = 0.2766948987373093
$ws.Cells.Item(5,10).Value = 0.2766948987373092
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.307106666666667
$ws.Cells.Item(5,14).Value = 3.92132
$ws.Cells.Item(5,15).Value = 0.01256263154946851
$ws.Cells.Item(5,16).Value = 0.01256263154946851
$ws.Cells.Item(5,17).Value = 14.81134456134667
$ws.Cells.Item(5,18).Value = 133.30210105212
$ws.Cells.Item(5,19).Value = 0.003476016064454316
$ws.Cells.Item(5,20).Value = 0.003476016064454316

# Row 6: FAPs -> FAPs  (Tnc/Egfr)
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tnc"
$ws.Cells.Item(6,3).Value = "Egfr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 11.331397
$ws.Cells.Item(6,8).Value = 33.994191
$ws.Cells.Item(6,9).Value = 0.2766948987373093
$ws.Cells.Item(6,10).Value = 0.2766948987373092
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 80.22623699999998
$ws.Cells.Item(6,14).Value = 240.678711
$ws.Cells.Item(6,15).Value = 0.77105616682495
$ws.Cells.Item(6,16).Value = 0.77105616682495
$ws.Cells.Item(6,17).Value = 909.0753412630888
$ws.Cells.Item(6,18).Value = 8181.6780713678
$ws.Cells.Item(6,19).Value = 0.2133473080004074
$ws.Cells.Item(6,20).Value = 0.2133473080004074

# Row 7: FAPs -> sCs  (Tnc/Egfr)
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tnc"
$ws.Cells.Item(7,3).Value = "Egfr"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 11.331397
$ws.Cells.Item(7,8).Value = 33.994191
$ws.Cells.Item(7,9).Value = 0.2766948987373093
$ws.Cells.Item(7,10).Value = 0.2766948987373092
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 22.51385866666667
$ws.Cells.Item(7,14).Value = 67.541576
$ws.Cells.Item(7,15).Value = 0.2163812016255815
$ws.Cells.Item(7,16).Value = 0.2163812016255815
$ws.Cells.Item(7,17).Value = 255.1134705538907
$ws.Cells.Item(7,18).Value = 2296.021234985016
$ws.Cells.Item(7,19).Value = 0.05987157467244757
$ws.Cells.Item(7,20).Value = 0.05987157467244756

# Row 8: sCs -> ECs  (Tnc/Egfr)
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Tnc"
$ws.Cells.Item(8,3).Value = "Egfr"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 28.17890933333333
$ws.Cells.Item(8,8).Value = 84.536728
$ws.Cells.Item(8,9).Value = 0.6880846610982287
$ws.Cells.Item(8,10).Value = 0.6880846610982286
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.307106666666667
$ws.Cells.Item(8,14).Value = 3.92132
$ws.Cells.Item(8,15).Value = 0.01256263154946851
$ws.Cells.Item(8,16).Value = 0.01256263154946851
$ws.Cells.Item(8,17).Value = 36.83284024899555
$ws.Cells.Item(8,18).Value = 331.49556224096
$ws.Cells.Item(8,19).Value = 0.008644154072217957
$ws.Cells.Item(8,20).Value = 0.008644154072217955

# Row 9: sCs -> FAPs  (Tnc/Egfr)
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Tnc"
$ws.Cells.Item(9,3).Value = "Egfr"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 28.17890933333333
$ws.Cells.Item(9,8).Value = 84.536728
$ws.Cells.Item(9,9).Value = 0.6880846610982287
$ws.Cells.Item(9,10).Value = 0.6880846610982286
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 80.22623699999998
$ws.Cells.Item(9,14).Value = 240.678711
$ws.Cells.Item(9,15).Value = 0.77105616682495
$ws.Cells.Item(9,16).Value = 0.77105616682495
$ws.Cells.Item(9,17).Value = 2260.687858577512
$ws.Cells.Item(9,18).Value = 20346.19072719761
$ws.Cells.Item(9,19).Value = 0.5305519212374451
$ws.Cells.Item(9,20).Value = 0.530551921237445

# Row 10: sCs -> sCs  (Tnc/Egfr)
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Tnc"
$ws.Cells.Item(10,3).Value = "Egfr"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 28.17890933333333
$ws.Cells.Item(10,8).Value = 84.536728
$ws.Cells.Item(10,9).Value = 0.6880846610982287
$ws.Cells.Item(10,10).Value = 0.6880846610982286
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 22.51385866666667
$ws.Cells.Item(10,14).Value = 67.541576
$ws.Cells.Item(10,15).Value = 0.2163812016255815
$ws.Cells.Item(10,16).Value = 0.2163812016255815
$ws.Cells.Item(10,17).Value = 634.415982111481
$ws.Cells.Item(10,18).Value = 5709.743839003329
$ws.Cells.Item(10,19).Value = 0.1488885857885657
$ws.Cells.Item(10,20).Value = 0.1488885857885657
